$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 29.174445
$ws.Range("B3").Value = 0.08317511
$ws.Range("B6").Value = 4999.9997
$ws.Range("B7").Value = 1.999999
$ws.Range("B8").Value = 1.9585686
$ws.Range("B9").Value = [double]"1.0022805e-06"
$ws.Range("B10").Value = 0.12468041
$ws.Range("B11").Value = 252.3337
$ws.Range("B12").Value = 6.2042942
$ws.Range("B13").Value = 0.85143389
$ws.Range("B14").Value = 76.045075
$ws.Range("B15").Value = 12.991144
$ws.Range("B16").Value = [double]"1.0002893e-06"
$ws.Range("B17").Value = [double]"1.0000088e-06"
$ws.Range("B18").Value = [double]"1.0000008e-06"
$ws.Range("B19").Value = [double]"1.0000109e-06"
$ws.Range("B20").Value = 10.111966
$ws.Range("B21").Value = 50.799708
$ws.Range("B22").Value = 0.97339172
$ws.Range("B23").Value = 160.0545
$ws.Range("B24").Value = 9.866707699999999
$ws.Range("B25").Value = 1.679041
$ws.Range("B26").Value = [double]"1.5287718e-14"
$ws.Range("B27").Value = 598.06091
$ws.Range("B28").Value = 0.2299311
$ws.Range("B29").Value = 163.78642
$ws.Range("B30").Value = 0.26884059
$ws.Range("B31").Value = 147.81167
$ws.Range("B32").Value = 0.0057417882
$ws.Range("B33").Value = 83.73015700000001
$ws.Range("B34").Value = 0.0010753622
$ws.Range("B35").Value = 121.41986
$ws.Range("B36").Value = [double]"1.8399146e-14"
$ws.Range("B37").Value = 172.96309
$ws.Range("B38").Value = 0.28585693
$ws.Range("B39").Value = 295.80542
$ws.Range("B40").Value = 0.17991137
$ws.Range("B41").Value = 945.69826
$ws.Range("B42").Value = 0.24771521
$ws.Range("B43").Value = 288.53057
$ws.Range("B44").Value = 0.059868212
$ws.Range("B45").Value = 294.25045
$ws.Range("B46").Value = 0.12866976
$ws.Range("B47").Value = 182.34025
$ws.Range("B48").Value = 0.16242264
$ws.Range("B49").Value = 235.17924
$ws.Range("B50").Value = 0.30013331
$ws.Range("B51").Value = 255.50264
$ws.Range("B52").Value = [double]"2.9063089e-15"
$ws.Range("B53").Value = 133.89504
$ws.Range("B54").Value = 0.23081236
$ws.Range("B55").Value = 176.08767
$ws.Range("B56").Value = 0.34090307
$ws.Range("B57").Value = 167.91428
$ws.Range("B58").Value = 0.72501381
$ws.Range("B59").Value = 287.0452
$ws.Range("B60").Value = 0.56048866
$ws.Range("B61").Value = 185.4575
$ws.Range("B62").Value = 0.052511486
$ws.Range("B63").Value = 221.49705
$ws.Range("B64").Value = [double]"1.2077726e-14"
$ws.Range("B65").Value = 493.96858
$ws.Range("B66").Value = 0.54254741
$ws.Range("B67").Value = 132.07398
$ws.Range("B68").Value = 0.39406761
$ws.Range("B69").Value = 231.88819
$ws.Range("B70").Value = [double]"4.9126975e-14"
$ws.Range("B71").Value = 179.35222
$ws.Range("B72").Value = 0.17099221
$ws.Range("B73").Value = 165.15027
$ws.Range("B74").Value = 0.42166825
$ws.Range("B75").Value = 267.84076
$ws.Range("B76").Value = 0.66013641
$ws.Range("B77").Value = 186.81564
$ws.Range("B78").Value = 0.53677823
$ws.Range("B79").Value = 4
$ws.Range("B80").Value = 0.27920627
$ws.Range("B81").Value = 0.77420349
$ws.Range("B82").Value = 0.4220426
$ws.Range("B83").Value = 0.8540212700000001
$ws.Range("B84").Value = 0.59759378
$ws.Range("B85").Value = 0.66286819
$ws.Range("B86").Value = 0.7302060500000001
$ws.Range("B87").Value = 0.10197067
$ws.Range("B88").Value = 0.79773147
$ws.Range("B89").Value = 0.36521309
$ws.Range("B90").Value = 0.56346992
$ws.Range("B91").Value = 0.8831129
$ws.Range("B92").Value = 0.6503373
$ws.Range("B93").Value = 0.69524423
$ws.Range("B94").Value = 0.39671789
$ws.Range("B95").Value = 0.55823273
$ws.Range("B96").Value = 0.08156432600000001
$ws.Range("B97").Value = 0.69915048
$ws.Range("B98").Value = 0.2670861
$ws.Range("B99").Value = 0.052599107
$ws.Range("B100").Value = 0.07144830100000001
$ws.Range("B101").Value = 0.36071917
$ws.Range("B102").Value = 0.27344852
$ws.Range("B103").Value = 0.68070928
$ws.Range("B104").Value = 0.42063908
$ws.Range("B105").Value = 0.077733918
$ws.Range("B106").Value = 0.37332627
$ws.Range("B107").Value = 393.99558
$ws.Range("B108").Value = 0.058572145
$ws.Range("B109").Value = 0.56007462
$ws.Range("B110").Value = 221.73097
$ws.Range("B111").Value = 0.044340509
$ws.Range("B112").Value = 0.077143589
$ws.Range("B113").Value = 2.9830517
$ws.Range("B114").Value = 0.21456842
$ws.Range("B115").Value = 0.13850618
$ws.Range("B116").Value = 0.062815053
$ws.Range("B117").Value = 0.12667908
$ws.Range("B118").Value = 0.085529562
$ws.Range("B119").Value = 0.090478449
$ws.Range("B120").Value = 0.010167636
$ws.Range("B121").Value = 0.14606702
$ws.Range("B122").Value = 0.19484385
$ws.Range("B123").Value = 0.044333208
$ws.Range("B124").Value = 0.028468924
$ws.Range("B125").Value = 0.045430496
$ws.Range("B126").Value = 0.37018086
$ws.Range("B127").Value = 0.08371921
$ws.Range("B128").Value = 0.041881024
$ws.Range("B129").Value = 0.061928936
$ws.Range("B130").Value = 0.00046406412
$ws.Range("B131").Value = 0.017480835
$ws.Range("B132").Value = 0.076185724
$ws.Range("B133").Value = [double]"2.5898602e-15"
$ws.Range("B134").Value = 0.06320778000000001
$ws.Range("B135").Value = 0.99987473
$ws.Range("B136").Value = 0.053240696
$ws.Range("B137").Value = 0.0031108321
$ws.Range("B138").Value = 0.13295025
$ws.Range("B139").Value = 0.055742829
$ws.Range("B140").Value = 0.10282592
$ws.Range("B141").Value = 0.06929489799999999
$ws.Range("B142").Value = 247.89902
$ws.Range("B143").Value = 0.25491456
$ws.Range("B144").Value = 0.11324645
$ws.Range("B145").Value = 0.015530967
$ws.Range("B146").Value = 470.07545
$ws.Range("B147").Value = 0.34122448
$ws.Range("B148").Value = [double]"2.1465414e-13"
$ws.Range("B149").Value = 0.3604926
$ws.Range("B150").Value = 105.7038
$ws.Range("B151").Value = 0.83925989
$ws.Range("B152").Value = 0.66302494
$ws.Range("B153").Value = 0.10549622
$ws.Range("B154").Value = 454.40396
$ws.Range("B155").Value = 0.32376543
$ws.Range("B156").Value = 0.037861226
$ws.Range("B157").Value = 0.18378648
$ws.Range("B158").Value = 227.82486
$ws.Range("B159").Value = 0.55737514
$ws.Range("B160").Value = 0.013638279
$ws.Range("B161").Value = 0.04906603
$ws.Range("B162").Value = 296.92859
$ws.Range("B163").Value = 0.13350682
$ws.Range("B164").Value = 0.42769883
$ws.Range("B165").Value = 0.40670903
$ws.Range("B166").Value = 0.86539727
$ws.Range("B167").Value = 10.616194
$ws.Range("B168").Value = 10.924129
$ws.Range("B169").Value = 1.4097157
$ws.Range("B170").Value = 12.82697
$ws.Range("B171").Value = 0.4588519
$ws.Range("B172").Value = 14.269684
